# The two observation records stored in rows 2 and 3 were reordered: the
# record that used to be row 3 is now row 2, and the record that used to
# be row 2 is now row 3. Apply that swap explicitly, column by column.
#
# String values are written with a leading "'" (quote-prefix) and the
# cell style is then reset to "Normal". This forces Excel's text-entry
# path so values that look like dates/numbers/times (e.g. "2022-12-07",
# "1") are stored as literal text - exactly like the source data - and
# keeps "blank but present" text cells blank-but-present instead of
# turning into genuinely empty (absent) cells. The quote-prefix style
# that briefly gets applied is removed again by the "Normal" reset, so
# no stray cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Set-Number($addr, $number) {
    $ws.Range($addr).Value = $number
}

# ---- Row 2 gets what used to be row 3's values ----
Set-Number "A2" 104986863
Set-Number "B2" 57193
Set-Number "E2" 206004
Set-Text   "F2" "Skogshare"
Set-Text   "G2" "Lepus timidus"
Set-Text   "H2" "Linnaeus, 1758"
Set-Text   "K2" ""
Set-Text   "L2" ""
Set-Text   "M2" "gående/springande"
Set-Text   "N2" ""
Set-Text   "P2" "Gläfse, Jörken, Dlr"
Set-Number "Q2" 537888.8853063835
Set-Number "R2" 6669232.05540918
Set-Number "S2" 50
Set-Text   "Y2" "2022-12-07"
Set-Text   "Z2" "13:30"
Set-Text   "AA2" "2022-12-07"
Set-Text   "AB2" "13:30"
$ws.Range("AC2").Value = $null
Set-Text   "AW2" "Lars Mattsson"
Set-Text   "AX2" "Lars Mattsson"

# ---- Row 3 gets what used to be row 2's values ----
Set-Number "A3" 106077283
Set-Number "B3" 56278
Set-Number "E3" 100011
Set-Text   "F3" "Kungsörn"
Set-Text   "G3" "Aquila chrysaetos"
Set-Text   "H3" "(Linnaeus, 1758)"
$ws.Range("K3").Value = $null
$ws.Range("L3").Value = $null
Set-Text   "M3" "förbiflygande"
$ws.Range("N3").Value = $null
Set-Text   "P3" "Styggtjärnsberget, Dlr"
Set-Number "Q3" 537773.3909779217
Set-Number "R3" 6668679.681769322
Set-Number "S3" 10
Set-Text   "Y3" "2022-03-09"
Set-Text   "Z3" "00:00"
Set-Text   "AA3" "2022-03-09"
Set-Text   "AB3" "00:00"
Set-Text   "AC3" "Mötte örn 1 och de kollade in varandra och tog några svängar tillsammans i två omgångar. Uppfattades av mig som uppvaktning, men ej att de var ett etablerat par. De skiljdes sedan och denna örn drog vidare åt väster. Åldern bedömdes t subad me"
Set-Text   "AW3" "Samuel Keith"
Set-Text   "AX3" "Samuel Keith"
